# Regenerate handoff report: new handoff GUID/hash and updated timestamps.
$wb = $excel.ActiveWorkbook

$oldGuid = "372f9ef2-7f78-4ca9-9355-9fba3350fe01"
$newGuid = "af3440f9-268c-4d65-899a-c191478fdc30"
$oldHash = "2eff4df5cfd8b0095f75fcbb54e786ccc24676ca"
$newHash = "9abf3384b910a4b42fb097567f9f61cbfb96a409"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"
$oldZh = "$oldGuid.$oldHash.zh-cn.xlf"
$newZh = "$newGuid.$newHash.zh-cn.xlf"
$oldDe = "$oldGuid.$oldHash.de-de.xlf"
$newDe = "$newGuid.$newHash.de-de.xlf"

# === Sheet "Overview" ===
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = "2016-46-12 14:46:40"

$ovAddrA2 = "https://github.com/OpenLocalizationTest/oltest/blob/3e9fdc6e1a5aa574eb2ac49128f9a38b00d4a386/e2e/$oldMd"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $ovAddrA2, "", "", $newMd)

# === Sheet "zh-cn" ===
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("D2").Value = $newZh
$wsZh.Range("E2").Value = "2016-03-12 14:46:36"

$zhAddrA2 = "https://github.com/OpenLocalizationTest/oltest/blob/3e9fdc6e1a5aa574eb2ac49128f9a38b00d4a386/e2e/$oldMd"
$zhAddrB2 = "https://github.com/OpenLocalizationTest/oltest/blob/3e9fdc6e1a5aa574eb2ac49128f9a38b00d4a386/e2e/$oldMd"
$zhAddrD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3aa366c8b1bfe97e307cc1430f7b54abeaf6f1c7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZh"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhAddrA2, "", "", $newMd)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $zhAddrB2, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhAddrD2, "", "", $newZh)

# === Sheet "de-de" ===
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("D2").Value = $newDe
$wsDe.Range("E2").Value = "2016-03-12 14:46:40"

$deAddrA2 = "https://github.com/OpenLocalizationTest/oltest/blob/3e9fdc6e1a5aa574eb2ac49128f9a38b00d4a386/e2e/$oldMd"
$deAddrB2 = "https://github.com/OpenLocalizationTest/oltest/blob/3e9fdc6e1a5aa574eb2ac49128f9a38b00d4a386/e2e/$oldMd"
$deAddrD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f285081af53e3f38ace91df5d07c741d12c608ec/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDe"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deAddrA2, "", "", $newMd)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $deAddrB2, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deAddrD2, "", "", $newDe)
